$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.892009437084198
$ws.Range("B1").Value = 1.482118844985962
$ws.Range("C1").Value = 2.942968606948853
$ws.Range("D1").Value = 3.99376106262207
$ws.Range("E1").Value = 1.200335741043091
